# --- Overview sheet ---
$wsOverview = $excel.ActiveWorkbook.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E1").ColumnWidth = 29.166666666666664
$wsOverview.Range("F1").ColumnWidth = 29.166666666666664

# --- zh-cn sheet ---
$wsZh = $excel.ActiveWorkbook.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("K2").Value = "2016-08-29 14:57:29"
$wsZh.Range("P2").Value = ""
$wsZh.Range("C1").ColumnWidth = 29.166666666666664
$wsZh.Range("P1").ColumnWidth = 12.833333333333332

# --- de-de sheet ---
$wsDe = $excel.ActiveWorkbook.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("K2").Value = "2016-08-29 14:57:36"
$wsDe.Range("P2").Value = ""
$wsDe.Range("C1").ColumnWidth = 29.166666666666664
$wsDe.Range("P1").ColumnWidth = 12.833333333333332
